$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Rumania (row 37)
$ws.Range("B37").Value = 17036
$ws.Range("C37").Value = 165
$ws.Range("D37").Value = 9930
$ws.Range("E37").Value = 5999

# Update Marruecos (row 59)
$ws.Range("B59").Value = 6930
$ws.Range("C59").Value = 60
$ws.Range("D59").Value = 3732
$ws.Range("E59").Value = 3006

# Update row 65
$ws.Range("E65").Value = 3858
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 25

# Move/update Etiopia: remove old row, insert new row with updated data
$ws.Rows("142:142").Delete()
$ws.Rows("135:135").Insert()
$ws.Range("A135").Value = "Etiopia"
$ws.Range("B135").Value = 352
$ws.Range("C135").Value = 35
$ws.Range("D135").Value = 116
$ws.Range("E135").Value = 231
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 5

# Update timestamp
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 12:35"
